$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the values from C58:F59 while preserving cell formatting/styles
$ws.Range("C58:F59").ClearContents()
